$wb = $excel.ActiveWorkbook

# 1) Grab the sheets by their current (pre-edit) names before renaming.
$infoSheet    = $wb.Worksheets.Item("INFO")
$templateSheet = $wb.Worksheets.Item("SubCalc_template")

# 2) Rename the template sheet (workbook.xml <sheet name="...">).
$templateSheet.Name = "subcalc-footprint-template"

# 3) Update the INFO sheet's long description cell (shared string text).
$infoSheet.Range("A1").Value = "The subcalc footprint template should be filled out as a flat file with entries in all columns for all rows. It can be kept in excel format or saved to a csv."

# 4) Move the template sheet's remembered selection to D39, then switch
#    back to the INFO sheet so it ends up as the active/selected tab again
#    (matches tabSelected="1" staying on INFO in the saved file).
[void]$templateSheet.Activate()
[void]$templateSheet.Range("D39").Select()

[void]$infoSheet.Activate()
[void]$infoSheet.Range("A1").Select()

Write-Host "Edits applied"
